$d = $word.ActiveDocument

$replacements = @(
    @("41×54=2214", "57×85=4845"),
    @("60×31=1860", "87×64=5568"),
    @("22×79=1738", "59×87=5133"),
    @("29×39=1131", "47×36=1692"),
    @("55×36=1980", "14×88=1232"),
    @("30×65=1950", "75×99=7425"),
    @("34×86=2924", "52×50=2600"),
    @("68×79=5372", "89×86=7654"),
    @("69×13=897", "49×65=3185"),
    @("59×95=5605", "41×85=3485"),
    @("90×96=8640", "65×60=3900"),
    @("83×53=4399", "77×96=7392"),
    @("13×52=676", "63×77=4851"),
    @("39×74=2886", "24×51=1224"),
    @("98×94=9212", "99×95=9405"),
    @("59×57=3363", "36×50=1800"),
    @("63×92=5796", "40×24=960"),
    @("78×28=2184", "45×25=1125"),
    @("64×39=2496", "73×84=6132"),
    @("50×66=3300", "18×34=612"),
    @("89×68=6052", "95×66=6270"),
    @("94×87=8178", "62×56=3472"),
    @("57×93=5301", "70×65=4550"),
    @("79×90=7110", "84×32=2688"),
    @("35×57=1995", "97×59=5723"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: not found -> $old"
    }
}
